# Swap the deck's theme colours: the slide master (ppt/theme/theme1.xml,
# currently the "Integral" theme) is re-coloured to match the stock
# "Office Theme" palette that the Notes Master (ppt/theme/theme2.xml)
# already uses.
#
# PowerPoint's object model addresses the 12 theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) through
# Master.Theme.ThemeColorScheme.Item(1..12).RGB, in that fixed order.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office" colour scheme (dk1..folHlink).
$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgb($officeColors[$i - 1])
}
